$d = $word.ActiveDocument

# Mapping of old equation text to new equation text, applied in document order.
# Each Find/Replace targets the whole document content; since every "old"
# string is unique at the time it is processed, ReplaceAll is safe here.
$replacements = @(
    @("91÷6=", "33÷5="),
    @("51÷2=", "58÷5="),
    @("42÷7=", "80÷8="),
    @("34÷8=", "98÷8="),
    @("32÷8=", "34÷2="),
    @("88÷9=", "47÷2="),
    @("61÷8=", "48÷6="),
    @("44÷7=", "96÷9="),
    @("49÷5=", "15÷3="),
    @("76÷7=", "58÷9="),
    @("53÷4=", "74÷3="),
    @("84÷5=", "73÷4="),
    @("56÷2=", "48÷5="),
    @("15÷7=", "56÷7="),
    @("19÷7=", "66÷4="),
    @("21÷9=", "14÷8="),
    @("24÷2=", "44÷7="),
    @("75÷7=", "96÷2="),
    @("66÷5=", "76÷2="),
    @("97÷8=", "67÷4="),
    @("72÷5=", "21÷2="),
    @("87÷3=", "82÷6="),
    @("49÷8=", "40÷8="),
    @("10÷5=", "89÷6="),
    @("25÷5=", "51÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
